$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing columns (U:AD) from rows 1 and 2
$ws.Range("U1:AD2").Clear()

# Reorder the Miller-index labels shown in row 2 (C2:J2)
$ws.Range("C2").Value = "[1, 1, 0]"
$ws.Range("D2").Value = "[2, 2, 2]"
$ws.Range("E2").Value = "[3, 1, 0]"
$ws.Range("F2").Value = "[3, 2, 1]"
$ws.Range("G2").Value = "[2, 1, 1]"
$ws.Range("H2").Value = "[4, 0, 0]"
$ws.Range("I2").Value = "[2, 2, 0]"
$ws.Range("J2").Value = "[2, 0, 0]"

# Append 4 new data rows (20-23), matching the formatting of the existing rows
$ws.Range("A19").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A20").Value = 18
$ws.Range("A21").Value = 19
$ws.Range("A22").Value = 20
$ws.Range("A23").Value = 21

$ws.Range("B20").Value = "HexGrid-90degTilt2.5degRes"
$ws.Range("B21").Value = "HexGrid-90degTilt5degRes"
$ws.Range("B22").Value = "HexGrid-90degTilt10degRes"
$ws.Range("B23").Value = "HexGrid-90degTilt15degRes"

$ws.Range("C20:T23").Value = 1
